$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this shifts existing rows 42-121 down to 43-122,
# carrying all of their data/formatting with them.
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new record.
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44477
$ws.Range("D42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 100112017
$ws.Range("G42").Value = "Apio"
$ws.Range("H42").Value = "Americana (o)"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 500
$ws.Range("K42").Value = 7000
$ws.Range("L42").Value = 7000
$ws.Range("M42").Value = 7000
$ws.Range("N42").Value = "$/docena de matas"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 1167
$ws.Range("Q42").Value = 6
$ws.Range("R42").Value = "Hortaliza"
